$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Row 5: Action renamed Change -> Update (set first so the new shared
# string "Update" is interned before the other new strings below)
$ws.Range("A5").Value = "Update"

# Row 4: Add / condition -- where clause now covers both TEST1 and TEST2,
# and records a successful run (Success=2, Failed=0)
$ws.Range("C4").Value = 'conditionnum in ["COG_TEST1","COG_TEST2"]'
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 0

# Row 5 (cont'd): condition now targets COG_TEST3; one successful run
$ws.Range("C5").Value = 'conditionnum in ["COG_TEST3"]'
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0

# Row 6: Action renamed Change -> Update; whereclause unchanged; one
# successful run recorded
$ws.Range("A6").Value = "Update"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0

# Row 7: Action renamed Change -> Update; object/whereclause unchanged;
# one successful run recorded, and the Error cell (F7) picks up the same
# wrap-text style already used by F4:F6
$ws.Range("A7").Value = "Update"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").WrapText = $true

# Update the active selection shown on the Input tab
$ws.Activate() | Out-Null
$ws.Range("B5").Select() | Out-Null
